# Auto-generated PowerShell COM-interop script to apply market-data refresh
# to the Mateus_Profits workbook (columns H-N on each Leve sheet).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 471.76923
$ws.Range("I28").Value = 136.08333
$ws.Range("K28").Value = 136.08333
$ws.Range("M28").Value = 348.91667
$ws.Range("H76").Value = 4197
$ws.Range("I76").Value = 4197
$ws.Range("K76").Value = 4197
$ws.Range("M76").Value = -3882
$ws.Range("H79").Value = 4197
$ws.Range("I79").Value = 4197
$ws.Range("K79").Value = 4197
$ws.Range("M79").Value = -3105
$ws.Range("H86").Value = 13915.889
$ws.Range("I86").Value = 14392
$ws.Range("J86").Value = 12249.5
$ws.Range("K86").Value = 14392
$ws.Range("L86").Value = 12249.5
$ws.Range("M86").Value = -13269
$ws.Range("N86").Value = -14495.5
$ws.Range("H88").Value = 4692.154
$ws.Range("I88").Value = 4149.6665
$ws.Range("J88").Value = 4854.9
$ws.Range("K88").Value = 4149.6665
$ws.Range("L88").Value = 4854.9
$ws.Range("M88").Value = -3743.6665
$ws.Range("N88").Value = -5666.9
$ws.Range("H89").Value = 13915.889
$ws.Range("I89").Value = 14392
$ws.Range("J89").Value = 12249.5
$ws.Range("K89").Value = 71960
$ws.Range("L89").Value = 61247.5
$ws.Range("M89").Value = -66344
$ws.Range("N89").Value = -72479.5
$ws.Range("H91").Value = 4692.154
$ws.Range("I91").Value = 4149.6665
$ws.Range("J91").Value = 4854.9
$ws.Range("K91").Value = 4149.6665
$ws.Range("L91").Value = 4854.9
$ws.Range("M91").Value = -2745.6665
$ws.Range("N91").Value = -7662.9
$ws.Range("H92").Value = 305.57144
$ws.Range("I92").Value = 253.625
$ws.Range("K92").Value = 253.625
$ws.Range("M92").Value = 994.375
$ws.Range("H100").Value = 1408
$ws.Range("J100").Value = 1564.3334
$ws.Range("L100").Value = 1564.3334
$ws.Range("N100").Value = -2646.3334
$ws.Range("H107").Value = 406.73685
$ws.Range("I107").Value = 398.73334
$ws.Range("J107").Value = 436.75
$ws.Range("K107").Value = 398.73334
$ws.Range("L107").Value = 436.75
$ws.Range("M107").Value = 1521.26666
$ws.Range("N107").Value = -4276.75
$ws.Range("H111").Value = 1524
$ws.Range("J111").Value = 798
$ws.Range("L111").Value = 2394
$ws.Range("N111").Value = -8528
$ws.Range("H115").Value = 472.16666
$ws.Range("I115").Value = 531.6
$ws.Range("J115").Value = 175
$ws.Range("K115").Value = 1594.8
$ws.Range("L115").Value = 525
$ws.Range("M115").Value = -27.80000000000018
$ws.Range("N115").Value = -3659
$ws.Range("H118").Value = 502.42856
$ws.Range("I118").Value = 502.42856
$ws.Range("K118").Value = 1507.28568
$ws.Range("M118").Value = 149.71432
$ws.Range("H125").Value = 733.9091
$ws.Range("I125").Value = 689.2857
$ws.Range("K125").Value = 6203.571300000001
$ws.Range("M125").Value = -3743.571300000001
$ws.Range("H137").Value = 3118.7778
$ws.Range("I137").Value = 2724.2144
$ws.Range("K137").Value = 8172.6432
$ws.Range("M137").Value = -5622.6432

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H34").Value = 0
$ws.Range("J34").Value = 0
$ws.Range("L34").Value = 0
$ws.Range("N34").ClearContents()
$ws.Range("H55").Value = 22300
$ws.Range("J55").Value = 19733.334
$ws.Range("L55").Value = 19733.334
$ws.Range("N55").Value = -20363.334
$ws.Range("H61").Value = 9138.294
$ws.Range("I61").Value = 9138.294
$ws.Range("K61").Value = 9138.294
$ws.Range("M61").Value = -8926.294
$ws.Range("H102").Value = 4691.9565
$ws.Range("I102").Value = 2662.6667
$ws.Range("K102").Value = 2662.6667
$ws.Range("M102").Value = -1040.6667
$ws.Range("H132").Value = 1384.6123
$ws.Range("I132").Value = 1384.6123
$ws.Range("K132").Value = 4153.8369
$ws.Range("M132").Value = -1623.8369
$ws.Range("H136").Value = 9138.294
$ws.Range("I136").Value = 9138.294
$ws.Range("K136").Value = 27414.882
$ws.Range("M136").Value = -24864.882

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H35").Value = 33207.2
$ws.Range("J35").Value = 53345.332
$ws.Range("L35").Value = 53345.332
$ws.Range("N35").Value = -53965.332
$ws.Range("H82").Value = 24804
$ws.Range("J82").Value = 37340
$ws.Range("L82").Value = 37340
$ws.Range("N82").Value = -38106
$ws.Range("H85").Value = 24804
$ws.Range("J85").Value = 37340
$ws.Range("L85").Value = 37340
$ws.Range("N85").Value = -39992
$ws.Range("H105").Value = 3303.5
$ws.Range("I105").Value = 3373.6667
$ws.Range("K105").Value = 3373.6667
$ws.Range("M105").Value = -1626.6667
$ws.Range("H134").Value = 6316.7144
$ws.Range("I134").Value = 6119.6665
$ws.Range("K134").Value = 18358.9995
$ws.Range("M134").Value = -15823.9995

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H41").Value = 21106.715
$ws.Range("J41").Value = 21624.5
$ws.Range("L41").Value = 21624.5
$ws.Range("N41").Value = -22480.5
$ws.Range("H74").Value = 39982.5
$ws.Range("J74").Value = 41643.332
$ws.Range("L74").Value = 41643.332
$ws.Range("N74").Value = -43391.332
$ws.Range("H77").Value = 39982.5
$ws.Range("J77").Value = 41643.332
$ws.Range("L77").Value = 124929.996
$ws.Range("N77").Value = -133665.996
$ws.Range("H100").Value = 120000
$ws.Range("J100").Value = 120000
$ws.Range("L100").Value = 120000
$ws.Range("N100").Value = -122164
$ws.Range("H107").Value = 903.7646999999999
$ws.Range("I107").Value = 319.53845
$ws.Range("J107").Value = 2802.5
$ws.Range("K107").Value = 319.53845
$ws.Range("L107").Value = 2802.5
$ws.Range("M107").Value = 1600.46155
$ws.Range("N107").Value = -6642.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H75").Value = 1052
$ws.Range("I75").Value = 863
$ws.Range("K75").Value = 2589
$ws.Range("M75").Value = -1591
$ws.Range("H76").Value = 843332.7
$ws.Range("H78").Value = 1052
$ws.Range("I78").Value = 863
$ws.Range("K78").Value = 7767
$ws.Range("M78").Value = -2775
$ws.Range("H79").Value = 843332.7
$ws.Range("H109").Value = 9650
$ws.Range("I109").Value = 0
$ws.Range("J109").Value = 9650
$ws.Range("K109").Value = 0
$ws.Range("L109").Value = 28950
$ws.Range("M109").ClearContents()
$ws.Range("N109").Value = -31030

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H27").Value = 24902.5
$ws.Range("I27").Value = 10000
$ws.Range("J27").Value = 39805
$ws.Range("K27").Value = 10000
$ws.Range("L27").Value = 39805
$ws.Range("M27").Value = -9834
$ws.Range("N27").Value = -40137
$ws.Range("H80").Value = 3675.5715
$ws.Range("I80").Value = 2421.6667
$ws.Range("J80").Value = 4616
$ws.Range("K80").Value = 2421.6667
$ws.Range("L80").Value = 4616
$ws.Range("M80").Value = -1423.6667
$ws.Range("N80").Value = -6612
$ws.Range("H83").Value = 3675.5715
$ws.Range("I83").Value = 2421.6667
$ws.Range("J83").Value = 4616
$ws.Range("K83").Value = 12108.3335
$ws.Range("L83").Value = 23080
$ws.Range("M83").Value = -7116.333500000001
$ws.Range("N83").Value = -33064
$ws.Range("H93").Value = 48263.832
$ws.Range("J93").Value = 48263.832
$ws.Range("L93").Value = 48263.832
$ws.Range("N93").Value = -52007.832

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H39").Value = 25000
$ws.Range("J39").Value = 25000
$ws.Range("L39").Value = 25000
$ws.Range("N39").Value = -25920

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H52").Value = 30000
$ws.Range("J52").Value = 0
$ws.Range("L52").Value = 0
$ws.Range("N52").ClearContents()
$ws.Range("H70").Value = 35000
$ws.Range("J70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("N70").ClearContents()
$ws.Range("H73").Value = 35000
$ws.Range("J73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("N73").ClearContents()
$ws.Range("H113").Value = 900.6429000000001
$ws.Range("I113").Value = 512.2273
$ws.Range("K113").Value = 1536.6819
$ws.Range("M113").Value = 633.3181
